$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Carte avant" block: B13 was 1001 (hex) -> 100 (hex); formula recalculates.
$ws.Range("B13").Value = 100
$ws.Range("C13").Formula = "=HEX2DEC(B13)"

# Keep C16's formula intact (B16/1008 unchanged) - re-assert so it stays a formula cell.
$ws.Range("C16").Formula = "=HEX2DEC(B16)"

# "Carte arrière" block: B18 was 1100 (hex) -> 110 (hex); formula recalculates.
$ws.Range("B18").Value = 110
$ws.Range("C18").Formula = "=HEX2DEC(B18)"

# Move the active selection, as it ended up after the debug edits.
[void]$ws.Range("D20").Select()
